$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.8484848484848485
$ws.Range("D3").Value = 0.9696969696969697
$ws.Range("H3").Value = 0.7776649746192893
$ws.Range("I3").Value = 0.05589712301121697
$ws.Range("J3").Value = 0.7575757575757576
$ws.Range("K3").Value = 125.2727272727273

$ws.Range("Q3").Value = 12
$ws.Range("R3").Value = 28
$ws.Range("S3").Value = 41
$ws.Range("T3").Value = 89
$ws.Range("U3").Value = 213

$ws.Range("V3").Value = 1925
$ws.Range("W3").Value = 1909
$ws.Range("X3").Value = 1896
$ws.Range("Y3").Value = 1848
$ws.Range("Z3").Value = 1724

$ws.Range("AF3").Value = 0.993805
$ws.Range("AG3").Value = 0.985545
$ws.Range("AH3").Value = 0.978833
$ws.Range("AI3").Value = 0.954053
$ws.Range("AJ3").Value = 0.890036
